# ajustes sanity semilla 6 en clases de portabilidad prepago y postpago
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Semilla 6")

# --- Core data edits ---
$ws.Range("G5").Value  = "tvcan1020Sem_6"
$ws.Range("F5").Value  = "app"
$ws.Range("E10").Value = "3046008607"
$ws.Range("C12").Value = "3046008593"
$ws.Range("E9").Value  = "3046008609"
$ws.Range("C14").Value = "3046008600"

# --- View / selection adjustments ---
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("E9").Select()
